# Updated symbol list on Thu Dec 15 13:32:09 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell while forcing Excel to
# store it as text (matching the source workbook's inline-string cells),
# then restore the cell's default ("Normal") style so no stray number
# format / style index is left behind.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Price (column D) updates for rows with unchanged coin identity ---
Set-TextValue "D2"  "264.44"
Set-TextValue "D3"  "22.73"
Set-TextValue "D4"  "6.225"
Set-TextValue "D5"  "0.06132"
Set-TextValue "D6"  "3.562"
Set-TextValue "D7"  "6.730"
Set-TextValue "D8"  "1.374"
Set-TextValue "D9"  "0.8136"
Set-TextValue "D10" "0.1593"
Set-TextValue "D11" "0.08208"
Set-TextValue "D12" "0.03398"
Set-TextValue "D13" "0.03178"
Set-TextValue "D14" "0.09238"
Set-TextValue "D15" "3.915"

Set-TextValue "D16" "0.001702"
$ws.Range("E16").Value = "15BitForexTokenBFBestin24h"

Set-TextValue "D17" "0.04847"
Set-TextValue "D18" "0.0006262"
Set-TextValue "D19" "0.006260"
Set-TextValue "D20" "0.001101"
Set-TextValue "D21" "0.003208"
Set-TextValue "D22" "0.0001506"
Set-TextValue "D23" "3.695"
Set-TextValue "D24" "2.253"
Set-TextValue "D25" "0.3388"
Set-TextValue "D26" "0.1272"
Set-TextValue "D27" "0.0002691"

Set-TextValue "D40" "0.04570"

# --- Rows 41-43 reordered: KickToken / BKEXToken / CEJI cycle ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1129"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003142"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003443"
$ws.Range("E43").Value = "42KickTokenKICK"

Set-TextValue "D44" "0.01075"
Set-TextValue "D45" "0.00006179"
Set-TextValue "D46" "0.00000000753"
Set-TextValue "D47" "0.7526"

Set-TextValue "D48" "0.07709"
$ws.Range("E48").Value = "47BOLOBOLO"

Set-TextValue "D49" "0.00002107"
Set-TextValue "D50" "0.01244"
